$d = $word.ActiveDocument

$texto2 = " Todos os módulos de nível 3 foram testados por bancadas de testes e as formas de ondas obtidas foram como o esperado. Os circuitos incluem: contador_mod10, contador_mod6, latch_sr, comb_logic (lógica combinacional do magnetron), mux, encoder, counter_freq e counter_non_recycling. "
$texto3 = " Todos os blocos digitais do nível 2 foram integrados, implementados e testados com bancadas de testes, e, com elas obtivemos formas de ondas de acordo com o esperado. Os blocos são: magnetron, entrada_timer_controle, decoder_7seg, timer_min_sec."

# --- Step 1: split the "2)" placeholder paragraph -----------------------
# Before: [... "2)" ...] [... "3)" ...]
# After split: [empty] ["2)"] [empty] [... "3)" ... untouched]
$rng = $d.Content
$ok1 = $rng.Find.Execute("2)", $true, $false, $false, $false, $false, $true, 1, $false, "^p2)^p", 2)

# --- Step 2: locate the new "2)" paragraph and the original "3)" one ----
$p2 = $null
$p3 = $null
foreach ($p in $d.Paragraphs) {
    if ($p2 -eq $null -and $p.Range.Text -eq "2)`r") {
        $p2 = $p
    } elseif ($p3 -eq $null -and $p.Range.Text -eq "3)`r") {
        $p3 = $p
    }
}

# --- Step 3: append the level-3 answer onto the "2)" paragraph as its own run
$end2 = $p2.Range.End - 1
$insertRng2 = $d.Range($end2, $end2)
$insertRng2.InsertAfter($texto2)

$p2 = $d.Paragraphs(12)
$answerStart2 = $p2.Range.End - 1 - $texto2.Length
$answerRng2 = $d.Range($answerStart2, $p2.Range.End - 1)
$answerRng2.Font.Bold = 1
$answerRng2.Font.Bold = 0

# --- Step 4: append the level-2 answer onto the original "3)" paragraph -
$end3 = $p3.Range.End - 1
$insertRng3 = $d.Range($end3, $end3)
$insertRng3.InsertAfter($texto3)

$answerStart3 = $p3.Range.End - 1 - $texto3.Length
$answerRng3 = $d.Range($answerStart3, $p3.Range.End - 1)
$answerRng3.Font.Bold = 1
$answerRng3.Font.Bold = 0

Write-Host "Done"
